# Applies the "23rd commit - Core alerts - Placeholder and template
# management TCs" change: inserts two new worksheets
# (CoreAlertsPlaceholderManagement, AlertsTemplateManagement) right
# after "CallCentreRoleManagement", populates them, and tweaks a couple
# of view-state properties (workbook firstSheet/activeTab, sheet tab
# selection/zoom/selection cells) to match the target workbook state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the two new worksheets right after CallCentreRoleManagement
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("CallCentreRoleManagement")

$placeholderSheet = $wb.Worksheets.Add($null, $anchor)
$placeholderSheet.Name = "CoreAlertsPlaceholderManagement"

$templateSheet = $wb.Worksheets.Add($null, $placeholderSheet)
$templateSheet.Name = "AlertsTemplateManagement"

# Reference cells (existing sheets) whose cell style exactly matches the
# new style combinations needed on the new sheets - copying formats from
# them keeps the workbook style table free of needless duplicate xf's.
$srcCoreUserMgmt   = $wb.Worksheets.Item("CoreUserManagement")
$srcCallCentre     = $wb.Worksheets.Item("CallCentreRoleManagement")
$srcDisposition    = $wb.Worksheets.Item("Disposition_master")
$srcUpdationDisp   = $wb.Worksheets.Item("Updation_of_Disposition")

# style "53": plain cell, border all sides, vertical=top wrapText=1
$styleWrapTop      = $srcCoreUserMgmt.Range("A2")
# style "54": header cell, bold/white font, blue fill, border all sides,
# horizontal=left vertical=top
$styleHeaderBox    = $srcCallCentre.Range("B1")
# style "55": plain cell, border all sides, vertical=top
$styleTop          = $srcDisposition.Range("C2")
# style "56": plain cell, border all sides, no special alignment
$stylePlainBorder  = $srcUpdationDisp.Range("C2")
# style "57": header cell, bold/white font, blue fill, border left+right
# only, horizontal=left vertical=top
$styleHeaderLR     = $srcUpdationDisp.Range("E1")
# style "58": plain cell, border all sides, applyFill (no pattern
# override visually - fillId 0), no special alignment
$stylePlainFill    = $srcUpdationDisp.Range("H2")

# ---------------------------------------------------------------------
# 2. CoreAlertsPlaceholderManagement content
# ---------------------------------------------------------------------
$placeholderSheet.Columns.Item(1).ColumnWidth = 42.36328125

$placeholderSheet.Range("A1").Value = "TestScenario"
$placeholderSheet.Range("B1").Value = "Run"
$placeholderSheet.Range("A2").Value = "CoreAlertsPlaceholderManagement"
$placeholderSheet.Range("B2").Value = "Yes"

$styleHeaderBox.Copy()
$placeholderSheet.Range("A1:B1").PasteSpecial(-4122)

$styleWrapTop.Copy()
$placeholderSheet.Range("A2").PasteSpecial(-4122)

$styleTop.Copy()
$placeholderSheet.Range("B2").PasteSpecial(-4122)

$placeholderSheet.Range("E12").Select()

# ---------------------------------------------------------------------
# 3. AlertsTemplateManagement content
# ---------------------------------------------------------------------
$templateSheet.Columns.Item(1).ColumnWidth = 24.6328125
$templateSheet.Columns.Item(10).ColumnWidth = 65.08984375

$headerValues = @("TestScenario","Run","NotificationType","CurrentStatus","Category","NotificationTypefortemplatecreation","InitialStatusfortemplatecreation","Categoryfortemplatecreation","TemplateNamefortemplatecreation","TemplateBody","EditInitialStatusfortemplatemodify","EditReason")
for ($i = 0; $i -lt $headerValues.Length; $i++) {
    $templateSheet.Cells.Item(1, $i + 1).Value = $headerValues[$i]
}

$dataValues = @("AlertsTemplateManagement","Yes","SMS","Active","Notification","EMAIL","Active","Promotions","Test","Automation testing improves software quality by identifying defects early. Selenium with Java helps in UI testing, ensuring reliability and efficiency in development.","Inactive","Test")
for ($i = 0; $i -lt $dataValues.Length; $i++) {
    $templateSheet.Cells.Item(2, $i + 1).Value = $dataValues[$i]
}

# Row 1 styles: A1:B1 = headerBox (border all sides); C1:L1 = headerLR
# (border left/right only)
$styleHeaderBox.Copy()
$templateSheet.Range("A1:B1").PasteSpecial(-4122)

$styleHeaderLR.Copy()
$templateSheet.Range("C1:L1").PasteSpecial(-4122)

# Row 2 styles
$styleWrapTop.Copy()
$templateSheet.Range("A2").PasteSpecial(-4122)

$styleTop.Copy()
$templateSheet.Range("B2").PasteSpecial(-4122)

$stylePlainFill.Copy()
$templateSheet.Range("C2,D2,F2").PasteSpecial(-4122)

$stylePlainBorder.Copy()
$templateSheet.Range("E2,G2,H2,I2,K2,L2").PasteSpecial(-4122)

# style "59" - new: white (theme background1) fill + wrapText, used for
# the long TemplateBody cell; no existing cell has this exact
# combination so it is built directly.
$templateSheet.Range("J2").Interior.ThemeColor = 2
$templateSheet.Range("J2").WrapText = $true
$templateSheet.Range("J2").Borders.LineStyle = 1
$templateSheet.Range("J2").Borders.Weight = 2
$templateSheet.Range("J2").Borders.ColorIndex = 64

$templateSheet.Rows.Item(2).RowHeight = 43.5

$templateSheet.Range("J17").Select()
$templateSheet.Application.ActiveWindow.Zoom = 80

# ---------------------------------------------------------------------
# 4. View-state tweaks
# ---------------------------------------------------------------------
# CallCentreRoleManagement loses its tabSelected flag and its prior
# selection is replaced.
$srcCallCentre.Range("G17").Select()

# AlertsTemplateManagement becomes the selected/active tab.
$templateSheet.Select()
$templateSheet.Range("J17").Select()

Write-Output "edit complete"
